# Auto-update draw results: append the 2025-12-22 Pick 3 draw as row 97.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 97

# Write the new row's values as string-typed formulas first (so Excel
# doesn't auto-coerce number-looking / date-looking text into numeric or
# date values), then flatten them to literal values via copy/paste-special
# so the cells end up holding plain text with no residual formula or
# number-format styling - matching how the rest of the sheet is stored.
$ws.Range("A$row").Formula = '="2025-12-22"'
$ws.Range("B$row").Formula = '="Pick 3"'
$ws.Range("C$row").Formula = '="251222"'
$ws.Range("D$row").Formula = '="9-5-3"'
$ws.Range("E$row").Formula = '="2025-12-22T21:40:59.810+04:00"'

$rangeAddr = "A" + $row + ":E" + $row
$dataRow = $ws.Range($rangeAddr)
$dataRow.Copy()
$dataRow.PasteSpecial(-4163)
